$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the header row (row 1) labels from the original
# code-style identifiers to human-readable titles.
$ws.Range("A1").Value = "Titulo Static"
$ws.Range("B1").Value = "Content Status"
$ws.Range("C1").Value = "Dato Static"
$ws.Range("D1").Value = "Titulo"
$ws.Range("E1").Value = "Fecha de Publicacion"

# Consolidate the separate per-column conditional formatting rules
# on C2:C7 / D2:D7 / E2:E7 (same yellow highlight) into a single
# rule covering C2:E7.
$ws.Range("D2:D7").FormatConditions.Delete()
$ws.Range("E2:E7").FormatConditions.Delete()

$fc = $ws.Range("C2:C7").FormatConditions.Item(1)
$fc.ModifyAppliesToRange($ws.Range("C2:E7"))
